$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: update last-updated timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 12:44 PM"

# --- 2. Industry Analysis sheet: update "1 Year" (F) column for rows 2-76 ---
$wsIndustry = $wb.Worksheets.Item("Industry Analysis")
$industryF = @{
  2 = 21.0016
  3 = -16.2396
  4 = 27.1317
  5 = -50.6494
  6 = 53.2813
  7 = -8.106199999999999
  8 = -9.552099999999999
  9 = 36.3756
  10 = -6.1314
  11 = 31.9081
  12 = -18.4955
  13 = 14.0155
  14 = -36.0718
  15 = -0.1622
  16 = 0.1459
  17 = -22.0012
  18 = 1.0561
  19 = -27.708
  20 = 47.7309
  21 = 12.0959
  22 = 95.1491
  23 = -50.2657
  24 = -13.3427
  25 = -9.9316
  26 = 5.8244
  27 = -32.7692
  28 = -24.8224
  29 = -18.4191
  30 = 25.8569
  31 = 58.4712
  32 = -3.3862
  33 = -6.3282
  34 = 27.7203
  35 = 4.4873
  36 = -4.9458
  37 = 3.6074
  38 = -23.3973
  39 = 8.7355
  40 = -5.8541
  41 = -8.3934
  42 = 20.3818
  43 = 14.3164
  44 = -12.6846
  45 = 28.4075
  46 = -1.1135
  47 = -37.1997
  48 = -29.8569
  49 = -27.5511
  50 = -49.7478
  51 = -51.8002
  52 = -38.5254
  53 = -12.4886
  54 = -5.0725
  55 = -17.7445
  56 = -26.636
  57 = -29.3361
  58 = -11.9574
  59 = -24.5687
  60 = -12.3
  61 = -10.9446
  62 = -17.1229
  63 = -9.5038
  64 = 54.2749
  65 = -43.4736
  66 = 13.2687
  67 = 12.7149
  68 = 24.8057
  69 = -17.0328
  70 = -6.8927
  71 = 13.6034
  72 = 3.9995
  73 = -16.226
  74 = -16.2448
  75 = 28.6924
  76 = 48.9752
}
foreach ($row in $industryF.Keys) {
  $wsIndustry.Cells.Item($row, 6).Value = $industryF[$row]
}

# --- 3. Stock List sheet: row 2 (CAPTRU-RE1) removed, rows shift up, new row appended ---
$wsStock = $wb.Worksheets.Item("Stock List")
$wsStock.Rows.Item(2).Delete()

$stockNew = @{
  2 = @("NIFTYCASE", 10.19, -0.5854, 0)
  3 = @("MOMENTUM30", 31.54, -0.6614, 0)
  4 = @("CANHLIFE", 118.46, 0.6286, 11253.7)
  5 = @("FLEXIADD", 10.64, -1.0233, 0)
  6 = @("MOENERGY", 36.3, -0.6568000000000001, 0)
  7 = @("MONIFTY100", 26.49, 0.3409, 0)
  8 = @("RUBICON", 652.65, -0.1453, 10752.4289)
  9 = @("CRAMC", 317.2, 2.3226, 6325.5208)
  10 = @("LGEINDIA", 1633.4, -0.946, 110870.6825)
  11 = @("TATACAP", 329.3, 0.1521, 139783.5374)
  12 = @("ELIQUID", 1004.85, 0.0408, 0)
  13 = @("WEWORK", 632.15, -2.4008, 8472.2803)
  14 = @("GROWWRLTY", 10.8, -0.4608, 0)
  15 = @("ADVANCE", 130.05, -5.2666, 836.0358)
  16 = @("OMFREIGHT", 88.90000000000001, -0.5926, 299.3747)
  17 = @("GLOTTIS", 72.73999999999999, -0.8587, 672.1394)
  18 = @("FABTECH", 237.72, 0.4734, 1056.6843)
  19 = @("PACEDIGITK", 218.85, 0.1327, 4723.9063)
  20 = @("JAINREC", 377.25, 1.2208, 13018.3623)
  21 = @("EPACKPEB", 301.45, 1.979, 3028.1254)
  22 = @("BMWVENTLTD", 69.25, 0, 600.5014)
  23 = @("STYL", 372.4, -0.8388, 6025.649)
  24 = @("JARO", 621.5, -1.4821, 1377.0134)
  25 = @("SOLARWORLD", 309.1, -0.6269, 2679.0517)
  26 = @("ARSSBL", 537.3, 4.7266, 3370.2277)
  27 = @("GANESHCP", 274.4, -2.7984, 1108.9312)
  28 = @("ATLANTAELE", 1003.05, -1.7436, 7713.116)
  29 = @("GKENERGY", 213.85, -0.7933, 4337.2472)
  30 = @("SAATVIKGL", 528.2, -1.3079, 6713.6863)
  31 = @("IVALUE", 281.45, -0.3364, 1506.8799)
  32 = @("VMSTMT", 70.03, -0.9056, 347.5674)
  33 = @("EUROPRATIK", 321.75, 0.8147, 3288.285)
  34 = @("SHRINGARMS", 229.31, -1.2616, 2211.284)
  35 = @("DEVX", 44.53, -0.3803, 401.605)
  36 = @("URBANCO", 148.9, -2.0459, 21380.5798)
  37 = @("SML100CASE", 10.36, -0.7663, 0)
  38 = @("AONEGOLD", 11.28, -0.2653, 0)
  39 = @("ELM250", 16.72, 0.1797, 0)
  40 = @("AMANTA", 122.52, 1.407, 475.7372)
  41 = @("CPEDU", 315.9, 1.8539, 574.7148999999999)
  42 = @("AHCL", 139.27, 3.1706, 740.2409)
  43 = @("STLNETWORK", 26.59, -0.412, 1297.3822)
  44 = @("VIKRAN", 98.05, -1.783, 2528.8166)
  45 = @("MANUFGBEES", 151.77, -1.011, 0)
  46 = @("MEIL", 461.15, -0.7319, 1274.1632)
  47 = @("GROWWNXT50", 70.29000000000001, -0.4109, 0)
  48 = @("SHREEJISPG", 270.05, -0.7899, 4399.6074)
  49 = @("GEMAROMA", 219.52, -0.876, 1146.7097)
  50 = @("PATELRMART", 219.31, -1.0646, 732.5069999999999)
  51 = @("VIKRAMSOLR", 322, -1.5892, 11647.2884)
  52 = @("LTGILTCASE", 29.67, 0.2365, 0)
  53 = @("REGAAL", 89.13, -0.8675, 915.5742)
  54 = @("BLUESTONE", 711.95, 0.1266, 10773.2539)
  55 = @("MOSILVER", 145.9, -1.5054, 0)
  56 = @("ALLTIME", 308.75, 2.66, 2022.5526)
  57 = @("JSWCEMENT", 134.98, -0.4793, 18402.6999)
  58 = @("SBILIQETF", 1012.94, 0.0296, 0)
  59 = @("HILINFRA", 77.23, -0.3998, 0)
  60 = @("GROWWPOWER", 10.28, -0.9634, 0)
  61 = @("LOTUSDEV", 177.82, 0.3669, 8690.485000000001)
  62 = @("MBEL", 450.2, -0.7714, 2572.8126)
  63 = @("LAXMIINDIA", 145.62, -1.1942, 761.1248000000001)
  64 = @("CPPLUS", 1322.1, -0.264, 15497.9053)
  65 = @("SHANTIGOLD", 241.57, -1.6409, 1741.6231)
  66 = @("MOGOLD", 119.65, -0.5403, 0)
  67 = @("BRIGHOTEL", 82.39, -0.9855, 3129.5229)
  68 = @("INDIQUBE", 212.64, -0.7561, 4465.6847)
  69 = @("EBGNG", 346.65, 3.2311, 3952.2092)
  70 = @("LIQGRWBEES", 1014.74, 0.0246, 0)
  71 = @("CHEMBONDCH", 153.35, -1.6987, 412.459)
  72 = @("GROWWNIFTY", 10.29, -0.3872, 0)
  73 = @("ANTHEM", 702.25, -0.1209, 39439.0658)
  74 = @("QUALITY30", 21.05, -0.8945, 0)
  75 = @("SMARTWORKS", 606.65, 2.0867, 6931.2448)
  76 = @("TRAVELFOOD", 1316.3, 0.1141, 17332.9705)
}
foreach ($row in $stockNew.Keys) {
  $vals = $stockNew[$row]
  $wsStock.Cells.Item($row, 2).Value = $vals[0]
  $wsStock.Cells.Item($row, 3).Value = $vals[0]
  $wsStock.Cells.Item($row, 4).Value = $vals[1]
  $wsStock.Cells.Item($row, 5).Value = $vals[2]
  $wsStock.Cells.Item($row, 8).Value = $vals[3]
}

# Row 76 (new TRAVELFOOD row) needs A/F/G set as well since it is a brand-new row
$wsStock.Cells.Item(76, 1).Value = "📋"
$wsStock.Cells.Item(76, 6).Value = "N/A"
$wsStock.Cells.Item(76, 7).Value = "N/A"